$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.5085426078324898
$ws.Cells.Item(2, 3).Value = 0.6312770643337786
$ws.Cells.Item(2, 4).Value = 0.5775225091701024
$ws.Cells.Item(2, 5).Value = 0.7599490174808455
$ws.Cells.Item(2, 6).Value = 0.5845391272395405
$ws.Cells.Item(2, 7).Value = 15

$ws.Cells.Item(3, 2).Value = 0.3604583747834282
$ws.Cells.Item(3, 3).Value = 0.5284892497187151
$ws.Cells.Item(3, 4).Value = 0.3749043735392947
$ws.Cells.Item(3, 5).Value = 0.6122943520393559
$ws.Cells.Item(3, 6).Value = 0.5136324539103049
$ws.Cells.Item(3, 7).Value = 14

$ws.Cells.Item(4, 2).Value = 0.2671644164584084
$ws.Cells.Item(4, 3).Value = 0.4711099050880848
$ws.Cells.Item(4, 4).Value = 0.3212841183425387
$ws.Cells.Item(4, 5).Value = 0.5668192995501641
$ws.Cells.Item(4, 6).Value = 0.5203199983321719
$ws.Cells.Item(4, 7).Value = 13

$ws.Cells.Item(5, 2).Value = 0.4282641160228051
$ws.Cells.Item(5, 3).Value = 0.548481490052795
$ws.Cells.Item(5, 4).Value = 0.4014706193385897
$ws.Cells.Item(5, 5).Value = 0.6336170920505457
$ws.Cells.Item(5, 6).Value = 0.4877336824715215
$ws.Cells.Item(5, 7).Value = 12

$ws.Cells.Item(6, 2).Value = 0.4297788858055521
$ws.Cells.Item(6, 3).Value = 0.5712726212628527
$ws.Cells.Item(6, 4).Value = 0.4264280183742889
$ws.Cells.Item(6, 5).Value = 0.6530145621456607
$ws.Cells.Item(6, 6).Value = 0.5156451691415619
$ws.Cells.Item(6, 7).Value = 11

$ws.Cells.Item(7, 2).Value = 0.3783814472866451
$ws.Cells.Item(7, 3).Value = 0.5371069045676157
$ws.Cells.Item(7, 4).Value = 0.3902711452809265
$ws.Cells.Item(7, 5).Value = 0.6247168520865485
$ws.Cells.Item(7, 6).Value = 0.5239790343878171
$ws.Cells.Item(7, 7).Value = 10

$ws.Cells.Item(8, 2).Value = 0.3326355726653664
$ws.Cells.Item(8, 3).Value = 0.4945997455591428
$ws.Cells.Item(8, 4).Value = 0.3233547881867544
$ws.Cells.Item(8, 5).Value = 0.5686429355815075
$ws.Cells.Item(8, 6).Value = 0.4891798334788348
$ws.Cells.Item(8, 7).Value = 9

$ws.Cells.Item(9, 2).Value = 0.3620304170102688
$ws.Cells.Item(9, 3).Value = 0.5407351962130211
$ws.Cells.Item(9, 4).Value = 0.3743586110673328
$ws.Cells.Item(9, 5).Value = 0.6118485197067431
$ws.Cells.Item(9, 6).Value = 0.5273032071390139
$ws.Cells.Item(9, 7).Value = 8

$ws.Cells.Item(10, 2).Value = 0.4970330926187166
$ws.Cells.Item(10, 3).Value = 0.5512696337016082
$ws.Cells.Item(10, 4).Value = 0.4018017306810739
$ws.Cells.Item(10, 5).Value = 0.6338783248235216
$ws.Cells.Item(10, 6).Value = 0.4249154521118755
$ws.Cells.Item(10, 7).Value = 7

$ws.Cells.Item(11, 2).Value = 0.4078060000532727
$ws.Cells.Item(11, 3).Value = 0.4615113944401893
$ws.Cells.Item(11, 4).Value = 0.259491323550082
$ws.Cells.Item(11, 5).Value = 0.5094029088551439
$ws.Cells.Item(11, 6).Value = 0.3343990248860763
$ws.Cells.Item(11, 7).Value = 6

